# Update column G ("K") values for rows 2-33 on the active sheet,
# matching the regenerated save_data (K instead of Strike#).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(4,7,1,3,4,7,6,6,3,1,7,5,5,5,6,4,6,6,12,4,6,8,5,4,5,1,6,6,3,1,2,5)

$row = 2
foreach ($val in $newValues) {
    $ws.Cells.Item($row, 7).Value = $val
    $row++
}
